# Weekly update: a new price-report row for "Acelga" (Macroferia Regional de
# Talca) is inserted at row 462, pushing the existing rows 462-482 down to
# 463-483 (all of their data stays intact, just shifted one row down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 462, shifting rows 462:482
# down to 463:483.
$ws.Rows(462).Insert()

# Populate the newly inserted row with this week's figures.
$ws.Range("A462").Value = 5
$ws.Range("B462").Value = 'Macroferia Regional de Talca'
$ws.Range("C462").Value = 'Maule'
$ws.Range("D462").Value = 45147
$ws.Range("E462").Value = 7
$ws.Range("F462").Value = 100112009
$ws.Range("G462").Value = 'Acelga'
$ws.Range("H462").Value = 'Sin especificar'
$ws.Range("I462").Value = 'Primera'
$ws.Range("J462").Value = 800
$ws.Range("K462").Value = 1500
$ws.Range("L462").Value = 1500
$ws.Range("M462").Value = 1500
$ws.Range("N462").Value = '$/docena de atados (4 kilos)'
$ws.Range("O462").Value = 'Región del Maule'
$ws.Range("P462").Value = 375
$ws.Range("Q462").Value = 4
$ws.Range("R462").Value = 'Hortaliza'
